$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Turn the existing totals sheet ("总计") into the new "2022-Q1"
#    per-fund holdings sheet (it keeps its original sheetId).
# ------------------------------------------------------------------
$fundSheet = $wb.Worksheets.Item("总计")

# Duplicate it first (while it still is a 4-column totals sheet) so the
# copy becomes the brand-new "总计" sheet placed right after it, with all
# the original sheet-level formatting (page margins, sheetPr, styles).
$fundSheet.Copy($null, $fundSheet)
$newTotalSheet = $wb.Worksheets.Item("总计 (2)")

# Rename the original sheet to its new name.
$fundSheet.Name = "2022-Q1"

# Pull in the 8-column header/style template used by every other quarter
# sheet (identical styling to what this sheet already used for its own
# 4-column header, just extended out to column H).
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$q4Sheet.Range("A1:H26").Copy($fundSheet.Range("A1:H26"))

# Header row
$fundSheet.Cells.Item(1,2).Value = '基金代码'
$fundSheet.Cells.Item(1,3).Value = '基金名称'
$fundSheet.Cells.Item(1,4).Value = '基金规模'
$fundSheet.Cells.Item(1,5).Value = '股票总仓位'
$fundSheet.Cells.Item(1,6).Value = '仓位占比'
$fundSheet.Cells.Item(1,7).Value = '持有市值(亿元)'
$fundSheet.Cells.Item(1,8).Value = '仓位排名'

# Columns B-G hold text (fund code / name / percentages kept as the
# original source strings, e.g. so leading zeros in fund codes survive).
$fundSheet.Range("B2:G26").NumberFormat = "@"

$fundData = @(
    @('970009','方正证券金立方一年持有期混合A','16.01','75.98','4.23','0.6772',6),
    @('970010','方正证券金立方一年持有期混合C','15.96','75.98','4.23','0.6751',6),
    @('001305','九泰天富改革新动力混合A','3.74','88.86','6.58','0.2461',4),
    @('013610','中信保诚前瞻优势混合','16.85','58.37','1.40','0.2359',7),
    @('001782','九泰久益灵活配置混合A','2.33','94.33','8.36','0.1948',3),
    @('217001','招商安泰混合','4.22','70.79','4.00','0.1688',2),
    @('001844','九泰久益灵活配置混合C','1.47','94.33','8.36','0.1229',3),
    @('011410','中信建投量化进取6个月持有期混合A','9.13','93.80','1.11','0.1013',4),
    @('460009','华泰柏瑞量化先行混合A','9.13','90.47','0.88','0.0803',8),
    @('002291','诺安安鑫灵活配置混合','2.19','81.55','3.43','0.0751',9),
    @('005632','鹏华量化先锋混合','3.10','92.91','1.64','0.0508',6),
    @('006401','先锋量化优选灵活配置混合A','0.86','93.48','5.40','0.0464',1),
    @('009912','九泰天富改革新动力混合C','0.59','88.86','6.58','0.0388',4),
    @('001017','泰达宏利改革动力量化策略灵活配置混合A','1.50','92.25','2.36','0.0354',3),
    @('229002','泰达宏利逆向策略混合','1.63','92.33','2.17','0.0354',3),
    @('009486','光大保德信瑞和混合A','0.53','91.77','5.69','0.0302',5),
    @('011411','中信建投量化进取6个月持有期混合C','2.15','93.80','1.11','0.0239',4),
    @('009487','光大保德信瑞和混合C','0.33','91.77','5.69','0.0188',5),
    @('006402','先锋量化优选灵活配置混合C','0.25','93.48','5.40','0.0135',1),
    @('008437','九泰行业优选灵活配置混合A','0.11','51.13','6.32','0.0070',3),
    @('004724','先锋聚元灵活配置混合A','0.09','94.56','5.45','0.0049',1),
    @('008438','九泰行业优选灵活配置混合C','0.06','51.13','6.32','0.0038',3),
    @('004725','先锋聚元灵活配置混合C','0.05','94.56','5.45','0.0027',1),
    @('010246','华泰柏瑞量化先行混合C','0.12','90.47','0.88','0.0011',8),
    @('003550','泰达宏利改革动力量化策略灵活配置混合C','0.01','92.25','2.36','0.0002',3),
)

$r = 2
foreach ($row in $fundData) {
    $fundSheet.Cells.Item($r, 1).Value = $r - 2
    $fundSheet.Cells.Item($r, 2).Value = $row[0]
    $fundSheet.Cells.Item($r, 3).Value = $row[1]
    $fundSheet.Cells.Item($r, 4).Value = $row[2]
    $fundSheet.Cells.Item($r, 5).Value = $row[3]
    $fundSheet.Cells.Item($r, 6).Value = $row[4]
    $fundSheet.Cells.Item($r, 7).Value = $row[5]
    $fundSheet.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ------------------------------------------------------------------
# 2) Populate the new "总计" sheet: the original 6 rows shift down one
#    row, and a new 2022-Q1 summary row is inserted at the top.
# ------------------------------------------------------------------
$newTotalSheet.Name = "总计"

# Shift the existing data rows (2-6) down to (3-7), keeping their styles.
$newTotalSheet.Range("A2:D6").Copy($newTotalSheet.Range("A3:D7"))

# Re-use row 3s per-column styling (identical to the old row 2) as the
# template for the freshly inserted row 2.
$newTotalSheet.Range("A3:D3").Copy($newTotalSheet.Range("A2:D2"))

$newTotalSheet.Cells.Item(2,1).Value = 0
$newTotalSheet.Cells.Item(2,2).Value = '2022-Q1'
$newTotalSheet.Cells.Item(2,3).Value = 25
$newTotalSheet.Cells.Item(2,4).Value = 2.89

$excel.ActiveSheet = $newTotalSheet
